$d = $word.ActiveDocument

# 1) Insert a new bullet "Auto Promotion" right after the paragraph that
#    ends with "...ordered by ED" (mirrors the formatting of that
#    paragraph, which already uses ListParagraph / numId 17 / 360 line
#    spacing / sz 24).
$p67 = $d.Paragraphs.Item(67)
if ($p67.Range.Text -notmatch "ordered by ED") {
    throw "unexpected paragraph 67 content: $($p67.Range.Text)"
}
$insertPoint = $d.Range($p67.Range.End, $p67.Range.End)
$insertPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(68)
$newPara.Range.Text = "Auto Promotion"

# 2) Remove the stray empty paragraph right after the "Queue DS" bullet.
$pQueue = $d.Paragraphs.Item(53)
if ($pQueue.Range.Text -notmatch "Queue DS") {
    throw "unexpected paragraph 53 content: $($pQueue.Range.Text)"
}
$d.Paragraphs.Item(54).Range.Delete()

# 3) Drop the empty placeholder paragraph right after the "Workload:"
#    heading (directly before the table that follows it).
$pWorkload = $d.Paragraphs.Item(27)
if ($pWorkload.Range.Text -notmatch "Workload:") {
    throw "unexpected paragraph 27 content: $($pWorkload.Range.Text)"
}
$d.Paragraphs.Item(28).Range.Delete()

# 4) Drop the empty heading-styled paragraph that sits just before the
#    "Workload:" heading paragraph.
$d.Paragraphs.Item(26).Range.Delete()

# 5) Give the "Action Class / Assigning Rovers..." table row an explicit
#    height (trHeight 2249 twips == 112.45 pt). Do this last: touching
#    Row.Height disturbs the live Paragraphs index cache afterwards.
$tbl = $d.Tables.Item(2)
$targetRow = $null
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    if ($tbl.Cell($r, 1).Range.Text -match "Action Class") {
        $targetRow = $tbl.Rows.Item($r)
        break
    }
}
if ($null -eq $targetRow) {
    throw "could not locate the Action Class row"
}
$targetRow.Height = 2249 / 20
